$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1, copying the
# existing header formatting (bold, border, centered) from AC1 ("Unnamed: 28").
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
